# Added progress bars for refresh status, display ack and download dataset
# (get datasets list only). This inserts 3 new "generic error" rows into the
# Error codes table on sheet1, updates the table/dimension to match, and
# adjusts the sheet view (scroll position / selection) left behind by the
# edit session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert the 3 new rows (in top-to-bottom order so row numbers below
#     each insertion point keep shifting correctly) -----------------------

# 1 new row just above the blank separator that used to sit at row 27
$ws.Range("A27:A27").EntireRow.Insert()

# 2 new rows (1 content + 1 blank separator) above old row 33
$ws.Range("A34:A35").EntireRow.Insert()

# 3 new rows (1 content + 2 blank separator) above old row 37
$ws.Range("A40:A42").EntireRow.Insert()

# --- Fill in the new content rows ----------------------------------------

$ws.Range("A27").Value = "Refresh status generic error"
$ws.Range("B27").Value = "ERR503"
$ws.Range("C27").Value = "yes"

$ws.Range("A34").Value = "Generic download report error (get datasets list)"
$ws.Range("B34").Value = "ERR704"
$ws.Range("C34").Value = "yes"

$ws.Range("A40").Value = "Display ack generic error"
$ws.Range("B40").Value = "ERR805"
$ws.Range("C40").Value = "yes"

# --- Update the existing "SOAP call" error description --------------------

$ws.Range("A7").Value = "Send message failed as result of SOAP call (probably due to username wrongly typed or to no connection problems)"

# --- Resize the table / autofilter to cover the new rows ------------------

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C72"))

# --- Restore the view the author left the sheet scrolled to ---------------

$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("A35").Select()
